# "Generate Report for handback" - refresh the handoff/handback timestamps
# for the first file's row (row 2) on each language sheet.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D2").Value = "2016-01-11 03:42:29"
$zhcn.Range("G2").Value = "2016-01-11 03:43:34"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D2").Value = "2016-01-11 03:42:48"
$dede.Range("G2").Value = "2016-01-11 03:44:01"
